$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings remain stored as text (matching source format)
$textCells = @("D5", "D6", "D7", "D8", "D9", "D10", "D12", "D16", "D18", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D30", "D31", "D32", "D33", "D34", "D36", "D37", "D38", "D39", "D40", "D43", "D44", "D45", "D46", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated values
$ws.Range("D2").Value = "98.205.04"
$ws.Range("E2").Value = "  -0.23%  "

$ws.Range("D3").Value = "3.413.63"
$ws.Range("E3").Value = "  -0.36%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").Value = "254.48"
$ws.Range("E5").Value = "  -0.77%  "

$ws.Range("D6").Value = "662.25"
$ws.Range("E6").Value = "  -2.84%  "

$ws.Range("D7").Value = "1.49"
$ws.Range("E7").Value = "  +1.83%  "

$ws.Range("D8").Value = "0.431"
$ws.Range("E8").Value = "  -1.16%  "

$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").Value = "1.04"
$ws.Range("E9").Value = "  -2.40%  "

$ws.Range("B10").Value = "USDC"
$ws.Range("C10").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D10").Value = "0.999"
$ws.Range("E10").Value = "  -0.03%  "

$ws.Range("D11").Value = "3.408.56"
$ws.Range("E11").Value = "  -0.44%  "

$ws.Range("D12").Value = "44.83"
$ws.Range("E12").Value = "  +6.49%  "

$ws.Range("E13").Value = "  -3.40%  "

$ws.Range("D14").Value = "97.985.09"
$ws.Range("E14").Value = "  -0.23%  "

$ws.Range("E15").Value = "  -4.96%  "

$ws.Range("D16").Value = "0.0000259"

$ws.Range("D17").Value = "4.043.15"
$ws.Range("E17").Value = "  -0.49%  "

$ws.Range("D18").Value = "9.20"
$ws.Range("E18").Value = "  +1.41%  "

$ws.Range("D19").Value = "3.384.91"
$ws.Range("E19").Value = "  -1.24%  "

$ws.Range("D20").Value = "18.36"
$ws.Range("E20").Value = "  +3.19%  "

$ws.Range("B21").Value = "Stellar"
$ws.Range("C21").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D21").Value = "0.528"
$ws.Range("E21").Value = "  -8.55%  "

$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "11.53"
$ws.Range("E22").Value = "  +2.99%  "

$ws.Range("D23").Value = "513.51"
$ws.Range("E23").Value = "  +0.28%  "

$ws.Range("D24").Value = "3.44"
$ws.Range("E24").Value = "  -0.70%  "

$ws.Range("D25").Value = "0.0000203"
$ws.Range("E25").Value = "  -2.14%  "

$ws.Range("D26").Value = "6.91"
$ws.Range("E26").Value = "  +4.28%  "

$ws.Range("D27").Value = "94.18"
$ws.Range("E27").Value = "  -7.19%  "

$ws.Range("E28").Value = "  -3.46%  "

$ws.Range("D29").Value = "3.588.71"
$ws.Range("E29").Value = "  -0.51%  "

$ws.Range("D30").Value = "12.03"
$ws.Range("E30").Value = "  +3.06%  "

$ws.Range("D31").Value = "0.144"
$ws.Range("E31").Value = "  -4.75%  "

$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "2.78"
$ws.Range("E32").Value = "  +5.73%  "

$ws.Range("B33").Value = "Dai"
$ws.Range("C33").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").Value = "  -0.19%  "

$ws.Range("D34").Value = "0.189"
$ws.Range("E34").Value = "  -4.98%  "

$ws.Range("E35").Value = "  -0.31%  "

$ws.Range("D36").Value = "0.567"
$ws.Range("E36").Value = "  -2.08%  "

$ws.Range("D37").Value = "29.25"
$ws.Range("E37").Value = "  -2.87%  "

$ws.Range("D38").Value = "8.01"
$ws.Range("E38").Value = "  -0.74%  "

$ws.Range("D39").Value = "1.50"
$ws.Range("E39").Value = "  -1.50%  "

$ws.Range("D40").Value = "526.53"
$ws.Range("E40").Value = "  -2.59%  "

$ws.Range("E41").Value = "  -0.68%  "

$ws.Range("E42").Value = "  -0.04%  "

$ws.Range("D43").Value = "0.866"
$ws.Range("E43").Value = "  -2.42%  "

$ws.Range("D44").Value = "24.42"
$ws.Range("E44").Value = "  -1.27%  "

$ws.Range("D45").Value = "1.75"
$ws.Range("E45").Value = "  -0.93%  "

$ws.Range("D46").Value = "0.0428"
$ws.Range("E46").Value = "  -2.83%  "

$ws.Range("E47").Value = "  -3.83%  "

$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").Value = "2.28"
$ws.Range("E48").Value = "  +6.80%  "

$ws.Range("B49").Value = "Filecoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D49").Value = "5.67"
$ws.Range("E49").Value = "  -4.50%  "

$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").Value = "8.64"
$ws.Range("E50").Value = "  -5.04%  "

$ws.Range("D51").Value = "56.10"
$ws.Range("E51").Value = "  +1.47%  "
